# Add 9 new data rows (50-58) to Sheet1, continuing the existing daily
# series (rows 2-49) with the same B:J values as the last existing row
# (row 49) and dates incrementing by one day each (45606 .. 45614).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B:J values replicated from row 49 (identical across all of rows 4-49)
$rowValues = @(116.4121952, 0.00170247, 0.008850780000000001, 0.06933635, 12792.90181321, 465.80531254, 0.24, 1.7904431, 485.38834923)

$startRow = 50
$startDate = 45606

for ($i = 0; $i -lt 9; $i++) {
    $r = $startRow + $i

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $startDate + $i

    # Match the date-column formatting used by the existing rows (style index 2:
    # bold font, thin border, centered/top alignment, custom date number format).
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dateCell.Font.Bold = $true
    $dateCell.HorizontalAlignment = -4108
    $dateCell.VerticalAlignment = -4160
    $dateCell.Borders.LineStyle = 1

    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $rowValues[$j]
    }
}
